# comment out unneeded crap in node
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "Simple EC2-vs-Lambda"
$ws2 = $wb.Worksheets.Item(2)   # "Complex EC2-vs-Lambda )"

# Remove sheet protection from the complex sheet.
$ws2.Unprotect()

# Two manual input tweaks that ripple through all the dependent formulas.
$ws2.Range("Q54").Value = 1
$ws2.Range("Q61").Value = 0.25

# Re-sort the EC2 instance price lookup table (M30:N57) alphabetically by
# instance name instead of by price.
$ws2.Sort.SortFields.Clear()
$ws2.Sort.SortFields.Add($ws2.Range("M30:M57"))
$ws2.Sort.SetRange($ws2.Range("M30:N57"))
$ws2.Sort.Apply()

# A handful of per-instance prices were also updated (these ripple into the
# VLOOKUP-driven rows 12-21 above).
$ws2.Range("N37").Value = 0.65          # g2.2xlarge
$ws2.Range("N38").Value = 0.3066        # hi1.4xlarge
$ws2.Range("N39").Value = 0.0184        # m1.large
$ws2.Range("N40").Value = 0.0115        # m1.medium
$ws2.Range("N42").Value = 0.0469        # m1.xlarge
$ws2.Range("N43").Value = 0.0961        # m2.4xlarge

# Make the complex sheet the active tab/sheet and leave a selection on it,
# matching the author's last-touched cell.
$ws2.Activate()
$ws2.Range("N46").Select()
